$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "45.767.62"
$ws.Range("E2").Value = "  +6.39%  "
$ws.Range("D3").Value = "2.399.03"
$ws.Range("E3").Value = "  +4.31%  "
$ws.Range("E4").Value = "  -0.86%  "
$ws.Range("D5").Value = "'113.37"
$ws.Range("E5").Value = "  +7.71%  "
$ws.Range("D6").Value = "'318.25"
$ws.Range("E6").Value = "  +2.67%  "
$ws.Range("E7").Value = "  +1.19%  "
$ws.Range("E8").Value = "  -0.44%  "
$ws.Range("D9").Value = "'0.628"
$ws.Range("E9").Value = "  +3.86%  "
$ws.Range("D10").Value = "'41.99"
$ws.Range("E10").Value = "  +6.25%  "
$ws.Range("D11").Value = "'0.0929"
$ws.Range("E11").Value = "  +2.45%  "
$ws.Range("E12").Value = "  +5.45%  "
$ws.Range("E13").Value = "  +2.36%  "
$ws.Range("D14").Value = "'1.01"
$ws.Range("E14").Value = "  +1.77%  "
$ws.Range("D15").Value = "'15.79"
$ws.Range("E15").Value = "  +3.52%  "
$ws.Range("D16").Value = "2.762.88"
$ws.Range("E16").Value = "  +4.02%  "
$ws.Range("D17").Value = "2.402.54"
$ws.Range("E17").Value = "  +4.23%  "
$ws.Range("D18").Value = "45.687.84"
$ws.Range("E18").Value = "  +6.57%  "
$ws.Range("E19").Value = "  +2.15%  "
$ws.Range("E20").Value = "  +3.31%  "
$ws.Range("D21").Value = "'13.35"
$ws.Range("E21").Value = "  -0.11%  "
$ws.Range("D22").Value = "'74.67"
$ws.Range("E22").Value = "  +1.54%  "
$ws.Range("D23").Value = "'3.53"
$ws.Range("E23").Value = "  +1.71%  "
$ws.Range("D24").Value = "'264.67"
$ws.Range("E24").Value = "  -0.89%  "
$ws.Range("E25").Value = "  +5.10%  "
$ws.Range("E26").Value = "  -0.49%  "
$ws.Range("B27").Value = "Cosmos"
$ws.Range("C27").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D27").Value = "'11.31"
$ws.Range("E27").Value = "  +3.36%  "
$ws.Range("B28").Value = "Filecoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D28").Value = "'7.57"
$ws.Range("E28").Value = "  +1.15%  "
$ws.Range("E29").Value = "  +2.65%  "
$ws.Range("D30").Value = "'39.09"
$ws.Range("E30").Value = "  +3.72%  "
$ws.Range("D31").Value = "'22.76"
$ws.Range("E31").Value = "  +2.44%  "
$ws.Range("D32").Value = "'0.0979"
$ws.Range("E32").Value = "  +13.30%  "
$ws.Range("D33").Value = "'172.63"
$ws.Range("E33").Value = "  +4.59%  "
$ws.Range("E34").Value = "  +4.61%  "
$ws.Range("E35").Value = "  +1.24%  "
$ws.Range("B36").Value = "Kaspa"
$ws.Range("C36").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D36").Value = "'0.118"
$ws.Range("E36").Value = "  +5.56%  "
$ws.Range("B37").Value = "RenderToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D37").Value = "'4.91"
$ws.Range("E37").Value = "  +6.04%  "
$ws.Range("D38").Value = "'4.14"
$ws.Range("E38").Value = "  +14.21%  "
$ws.Range("E39").Value = "  +7.95%  "
$ws.Range("E40").Value = "  +3.02%  "
$ws.Range("D41").Value = "'1.78"
$ws.Range("E41").Value = "  +14.32%  "
$ws.Range("D42").Value = "'102.73"
$ws.Range("E42").Value = "  -3.87%  "
$ws.Range("E43").Value = "  +5.47%  "
$ws.Range("D44").Value = "'13.54"
$ws.Range("E44").Value = "  +10.04%  "
$ws.Range("D45").Value = "'71.78"
$ws.Range("E45").Value = "  +0.34%  "
$ws.Range("D46").Value = "'87.68"
$ws.Range("E46").Value = "  +14.95%  "
$ws.Range("E47").Value = "  -0.53%  "
$ws.Range("D48").Value = "'115.36"
$ws.Range("E48").Value = "  +3.70%  "
$ws.Range("D49").Value = "'5.65"
$ws.Range("E49").Value = "  +8.91%  "
$ws.Range("E50").Value = "  +6.68%  "
$ws.Range("D51").Value = "1.660.62"
$ws.Range("E51").Value = "  -3.42%  "
